# Update "F" column (visitor/view count) figures across the four
# worksheets of the 广州-漫展信息 workbook, per the regenerated site data.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value2  = 2465
$ws1.Cells.Item(3, 6).Value2  = 749
$ws1.Cells.Item(4, 6).Value2  = 248
$ws1.Cells.Item(5, 6).Value2  = 420
$ws1.Cells.Item(6, 6).Value2  = 709
$ws1.Cells.Item(8, 6).Value2  = 912
$ws1.Cells.Item(9, 6).Value2  = 573
$ws1.Cells.Item(10, 6).Value2 = 947
$ws1.Cells.Item(14, 6).Value2 = 73
$ws1.Cells.Item(16, 6).Value2 = 1102
$ws1.Cells.Item(17, 6).Value2 = 24400
$ws1.Cells.Item(18, 6).Value2 = 2314
$ws1.Cells.Item(19, 6).Value2 = 148
$ws1.Cells.Item(22, 6).Value2 = 73
$ws1.Cells.Item(23, 6).Value2 = 361
$ws1.Cells.Item(25, 6).Value2 = 80
$ws1.Cells.Item(26, 6).Value2 = 239
$ws1.Cells.Item(28, 6).Value2 = 73
$ws1.Cells.Item(30, 6).Value2 = 359
$ws1.Cells.Item(33, 6).Value2 = 189

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(7, 6).Value2  = 271
$ws2.Cells.Item(8, 6).Value2  = 109
$ws2.Cells.Item(9, 6).Value2  = 28
$ws2.Cells.Item(10, 6).Value2 = 247
$ws2.Cells.Item(11, 6).Value2 = 3641
$ws2.Cells.Item(13, 6).Value2 = 155
$ws2.Cells.Item(15, 6).Value2 = 18
$ws2.Cells.Item(19, 6).Value2 = 137

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(4, 6).Value2 = 797

# Sheet 4: 全部类型 (All types, combined listing)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(4, 6).Value2  = 2465
$ws4.Cells.Item(5, 6).Value2  = 797
$ws4.Cells.Item(6, 6).Value2  = 749
$ws4.Cells.Item(7, 6).Value2  = 248
$ws4.Cells.Item(8, 6).Value2  = 420
$ws4.Cells.Item(9, 6).Value2  = 709
$ws4.Cells.Item(14, 6).Value2 = 271
$ws4.Cells.Item(16, 6).Value2 = 912
$ws4.Cells.Item(17, 6).Value2 = 573
$ws4.Cells.Item(18, 6).Value2 = 947
$ws4.Cells.Item(21, 6).Value2 = 73
$ws4.Cells.Item(23, 6).Value2 = 1102
$ws4.Cells.Item(24, 6).Value2 = 24400
$ws4.Cells.Item(25, 6).Value2 = 28
$ws4.Cells.Item(26, 6).Value2 = 247
$ws4.Cells.Item(28, 6).Value2 = 155
$ws4.Cells.Item(30, 6).Value2 = 2314
$ws4.Cells.Item(31, 6).Value2 = 148
$ws4.Cells.Item(32, 6).Value2 = 18
$ws4.Cells.Item(36, 6).Value2 = 361
$ws4.Cells.Item(38, 6).Value2 = 239
$ws4.Cells.Item(41, 6).Value2 = 73
$ws4.Cells.Item(43, 6).Value2 = 137
$ws4.Cells.Item(47, 6).Value2 = 189
